# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the row that used to
# describe 3f333ecc-....md (status "Handed back: in sync with en-US")
# is now reporting 37e5cb7e-....md as "Ready for handoff" with a new
# handoff timestamp / error detail, and the two tracked files swap row
# order (3f333ecc now listed first, 37e5cb7e second) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 / Row 3 file identity swap (A/B columns)
$ov.Range("A2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$ov.Range("B2").Value = "e2e\3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$ov.Range("A3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$ov.Range("B3").Value = "e2e\37e5cb7e-861c-40ec-816c-c1383e08f148.md"

# Status / datetime for row 3 (now 37e5cb7e) -> "Ready for handoff"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-12 06:59:31"

# Hyperlinks on column B: addresses (r:id targets) stay the same, only
# the displayed text needs to follow the swapped file names.
$addrOv2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$addrOv3 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/3f333ecc-78b8-442d-9710-3b4ca4700805.md"

$ov.Range("B2:B3").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $addrOv2, "", "", "e2e\3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$ov.Hyperlinks.Add($ov.Range("B3"), $addrOv3, "", "", "e2e\37e5cb7e-861c-40ec-816c-c1383e08f148.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 now describes 3f333ecc, row 3 now describes 37e5cb7e.
$zh.Range("A2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$zh.Range("G2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.zh-cn.xlf"
$zh.Range("I2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$zh.Range("J2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.zh-cn.xlf"

$zh.Range("A3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-12 06:59:25"
$zh.Range("I3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$zh.Range("J3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/9345487ed68f66fff68badff9baf7b1cce089f87/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md."

$zhAddrA2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$zhAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8495ab1a14f008c55a4d60a7023c0d31f6732ead/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$zhAddrA3 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$zhAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8495ab1a14f008c55a4d60a7023c0d31f6732ead/e2e/3f333ecc-78b8-442d-9710-3b4ca4700805.md"

$zh.Range("A2:A3,I2:I3").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhAddrA2, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$zh.Hyperlinks.Add($zh.Range("I2"), $zhAddrI2, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$zh.Hyperlinks.Add($zh.Range("A3"), $zhAddrA3, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $zhAddrI3, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")

# "Error Detail" column got a lot wider to fit the new message.
$zh.Columns.Item(16).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$de.Range("G2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.de-de.xlf"
$de.Range("I2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$de.Range("J2").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.de-de.xlf"

$de.Range("A3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.de-de.xlf"
$de.Range("H3").Value = "2016-08-12 06:59:31"
$de.Range("I3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$de.Range("J3").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/9345487ed68f66fff68badff9baf7b1cce089f87/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md."

$deAddrA2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$deAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dfe88ca02821bf128fd8d1982c856dde5d4c0bad/e2e/37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$deAddrA3 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/025a10c9f8da6f84cd14c2429a7480d566a85698/e2e/3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$deAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dfe88ca02821bf128fd8d1982c856dde5d4c0bad/e2e/3f333ecc-78b8-442d-9710-3b4ca4700805.md"

$de.Range("A2:A3,I2:I3").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deAddrA2, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$de.Hyperlinks.Add($de.Range("I2"), $deAddrI2, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$de.Hyperlinks.Add($de.Range("A3"), $deAddrA3, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$de.Hyperlinks.Add($de.Range("I3"), $deAddrI3, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")

$de.Columns.Item(16).ColumnWidth = 39.1

Write-Host "Report regenerated for handoff."
